$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = $ws.Range("D223").NumberFormat()

$row170 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44463, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 40, 5000, 5000, 5000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2500, 2, 'Hortaliza')
$row171 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44196, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 30, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row172 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44301, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 90, 5000, 6000, 5556, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2778, 2, 'Hortaliza')
$row173 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44301, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 20, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Región Metropolitana', 2000, 2, 'Hortaliza')
$row174 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44301, 9, 100112040, 'Cilantro', 'Sin especificar', 'Segunda', 40, 3500, 3500, 3500, '$/docena de atados (2 kilos)', 'Región Metropolitana', 1750, 2, 'Hortaliza')
$row175 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44251, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 75, 6500, 7000, 6800, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 3400, 2, 'Hortaliza')
$row176 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44243, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 65, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row177 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44252, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 55, 6500, 7000, 6682, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 3341, 2, 'Hortaliza')
$row178 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44166, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 155, 2500, 3000, 2774, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 1387, 2, 'Hortaliza')
$row179 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44168, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 155, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row180 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44369, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 40, 4000, 5000, 4500, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2250, 2, 'Hortaliza')
$row181 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44433, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 65, 5000, 5000, 5000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2500, 2, 'Hortaliza')
$row182 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44221, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 110, 3000, 3000, 3000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 1500, 2, 'Hortaliza')
$row183 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44316, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 30, 5000, 5000, 5000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2500, 2, 'Hortaliza')
$row184 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44279, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 65, 6000, 6000, 6000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 3000, 2, 'Hortaliza')
$row185 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44397, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 20, 8000, 8000, 8000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 4000, 2, 'Hortaliza')
$row186 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44397, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 30, 6600, 6600, 6600, '$/docena de atados (2 kilos)', 'Región Metropolitana', 3300, 2, 'Hortaliza')
$row187 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44363, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 45, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row188 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44277, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 100, 5000, 6000, 5550, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2775, 2, 'Hortaliza')
$row189 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44291, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 65, 7000, 7000, 7000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 3500, 2, 'Hortaliza')
$row190 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44273, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 40, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row191 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44438, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 40, 5000, 5000, 5000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2500, 2, 'Hortaliza')
$row192 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44438, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 20, 5000, 5000, 5000, '$/docena de atados (2 kilos)', 'Región Metropolitana', 2500, 2, 'Hortaliza')
$row193 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44372, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 30, 6000, 6000, 6000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 3000, 2, 'Hortaliza')
$row194 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44372, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 40, 6000, 6000, 6000, '$/docena de atados (2 kilos)', 'Región del Maule', 3000, 2, 'Hortaliza')
$row195 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44286, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 50, 6000, 6000, 6000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 3000, 2, 'Hortaliza')
$row196 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44209, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 140, 3000, 3500, 3286, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 1643, 2, 'Hortaliza')
$row197 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44356, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 20, 3000, 3000, 3000, '$/docena de atados (1 kilo)', 'Región Metropolitana', 3000, 1, 'Hortaliza')
$row198 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44356, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 20, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row199 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44160, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 50, 3000, 4000, 3500, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 1750, 2, 'Hortaliza')
$row200 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44351, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 95, 4000, 5000, 4632, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2316, 2, 'Hortaliza')
$row201 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44365, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 110, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row202 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44306, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 65, 6000, 6000, 6000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 3000, 2, 'Hortaliza')
$row203 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44215, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 100, 3000, 3000, 3000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 1500, 2, 'Hortaliza')
$row204 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44175, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 150, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row205 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44461, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 25, 6000, 6000, 6000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 3000, 2, 'Hortaliza')
$row206 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44461, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 40, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Región Metropolitana', 2000, 2, 'Hortaliza')
$row207 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44357, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 50, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row208 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44203, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 200, 3000, 3000, 3000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 1500, 2, 'Hortaliza')
$row209 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44162, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 50, 3000, 3000, 3000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 1500, 2, 'Hortaliza')
$row210 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44410, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 30, 5000, 5000, 5000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2500, 2, 'Hortaliza')
$row211 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44410, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 70, 4300, 4300, 4300, '$/docena de atados (2 kilos)', 'Región Metropolitana', 2150, 2, 'Hortaliza')
$row212 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44410, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 40, 4300, 4300, 4300, '$/docena de atados (2 kilos)', 'Región del Maule', 2150, 2, 'Hortaliza')
$row213 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44411, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 10, 5000, 5000, 5000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2500, 2, 'Hortaliza')
$row214 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44411, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 20, 5000, 5000, 5000, '$/docena de atados (2 kilos)', 'Región Metropolitana', 2500, 2, 'Hortaliza')
$row215 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44257, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 40, 6000, 7000, 6500, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 3250, 2, 'Hortaliza')
$row216 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44176, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 50, 5000, 5000, 5000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2500, 2, 'Hortaliza')
$row217 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44176, 9, 100112040, 'Cilantro', 'Sin especificar', 'Segunda', 10, 3000, 3000, 3000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 1500, 2, 'Hortaliza')
$row218 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44239, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 110, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row219 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44292, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 55, 8000, 8000, 8000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 4000, 2, 'Hortaliza')
$row220 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44358, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 50, 4000, 4000, 4000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2000, 2, 'Hortaliza')
$row221 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44211, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 95, 3000, 3000, 3000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 1500, 2, 'Hortaliza')
$row222 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44425, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 10, 5000, 5000, 5000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 2500, 2, 'Hortaliza')
$row223 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44425, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 40, 4300, 5000, 4650, '$/docena de atados (2 kilos)', 'Región Metropolitana', 2325, 2, 'Hortaliza')
$row224 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44323, 9, 100112040, 'Cilantro', 'Sin especificar', 'Primera', 35, 6000, 6000, 6000, '$/docena de atados (2 kilos)', 'Provincia de Cautín', 3000, 2, 'Hortaliza')

$rows = @($row170, $row171, $row172, $row173, $row174, $row175, $row176, $row177, $row178, $row179, $row180, $row181, $row182, $row183, $row184, $row185, $row186, $row187, $row188, $row189, $row190, $row191, $row192, $row193, $row194, $row195, $row196, $row197, $row198, $row199, $row200, $row201, $row202, $row203, $row204, $row205, $row206, $row207, $row208, $row209, $row210, $row211, $row212, $row213, $row214, $row215, $row216, $row217, $row218, $row219, $row220, $row221, $row222, $row223, $row224)
$rowNums = @(170, 171, 172, 173, 174, 175, 176, 177, 178, 179, 180, 181, 182, 183, 184, 185, 186, 187, 188, 189, 190, 191, 192, 193, 194, 195, 196, 197, 198, 199, 200, 201, 202, 203, 204, 205, 206, 207, 208, 209, 210, 211, 212, 213, 214, 215, 216, 217, 218, 219, 220, 221, 222, 223, 224)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowVals = $rows[$r]
    $rowNum = $rowNums[$r]
    $arr = New-Object "object[,]" 1,18
    for ($i = 0; $i -lt 18; $i++) {
        $arr[0, $i] = $rowVals[$i]
    }
    $rng = $ws.Range("A" + $rowNum + ":R" + $rowNum)
    $rng.Value = $arr
    $ws.Range("D" + $rowNum).NumberFormat = $dateFmt
}

$null = $ws.Range("A1").Select()